$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1457.9077  # H15: 1658.25 -> 1457.9077
$ws.Cells.Item(15, 9).Value = 1457.9077  # I15: 1658.25 -> 1457.9077
$ws.Cells.Item(15, 11).Value = 4373.7231  # K15: 4974.75 -> 4373.7231
$ws.Cells.Item(15, 13).Value = -4204.7231  # M15: -4805.75 -> -4204.7231

$ws.Cells.Item(17, 8).Value = 78093.08  # H17: 78100.766 -> 78093.08
$ws.Cells.Item(17, 10).Value = 78093.08  # J17: 78100.766 -> 78093.08
$ws.Cells.Item(17, 12).Value = 234279.24  # L17: 234302.298 -> 234279.24
$ws.Cells.Item(17, 14).Value = -234615.24  # N17: -234638.298 -> -234615.24

$ws.Cells.Item(33, 8).Value = 267.33334  # H33: 137.3125 -> 267.33334
$ws.Cells.Item(33, 9).Value = 288.9091  # I33: 156.84616 -> 288.9091
$ws.Cells.Item(33, 10).Value = 30  # J33: 52.666668 -> 30
$ws.Cells.Item(33, 11).Value = 288.9091  # K33: 156.84616 -> 288.9091
$ws.Cells.Item(33, 12).Value = 30  # L33: 52.666668 -> 30
$ws.Cells.Item(33, 13).Value = -59.90910000000002  # M33: 72.15384 -> -59.90910000000002
$ws.Cells.Item(33, 14).Value = -488  # N33: -510.666668 -> -488

$ws.Cells.Item(80, 8).Value = 4538.407  # H80: 4729.115 -> 4538.407
$ws.Cells.Item(80, 9).Value = 520.6667  # I80: 586.0909 -> 520.6667
$ws.Cells.Item(80, 10).Value = 7752.6  # J80: 7767.3335 -> 7752.6
$ws.Cells.Item(80, 11).Value = 1562.0001  # K80: 1758.2727 -> 1562.0001
$ws.Cells.Item(80, 12).Value = 23257.8  # L80: 23302.0005 -> 23257.8
$ws.Cells.Item(80, 13).Value = -564.0001  # M80: -760.2727 -> -564.0001
$ws.Cells.Item(80, 14).Value = -25253.8  # N80: -25298.0005 -> -25253.8

$ws.Cells.Item(83, 8).Value = 4538.407  # H83: 4729.115 -> 4538.407
$ws.Cells.Item(83, 9).Value = 520.6667  # I83: 586.0909 -> 520.6667
$ws.Cells.Item(83, 10).Value = 7752.6  # J83: 7767.3335 -> 7752.6
$ws.Cells.Item(83, 11).Value = 4686.0003  # K83: 5274.8181 -> 4686.0003
$ws.Cells.Item(83, 12).Value = 69773.40000000001  # L83: 69906.0015 -> 69773.40000000001
$ws.Cells.Item(83, 13).Value = 305.9997000000003  # M83: -282.8181000000004 -> 305.9997000000003
$ws.Cells.Item(83, 14).Value = -79757.40000000001  # N83: -79890.0015 -> -79757.40000000001

$ws.Cells.Item(127, 8).Value = 1154.25  # H127: 1205.0435 -> 1154.25
$ws.Cells.Item(127, 9).Value = 508.66666  # I127: 601.1667 -> 508.66666
$ws.Cells.Item(127, 10).Value = 1303.2307  # J127: 1295.625 -> 1303.2307
$ws.Cells.Item(127, 11).Value = 1525.99998  # K127: 1803.5001 -> 1525.99998
$ws.Cells.Item(127, 12).Value = 3909.6921  # L127: 3886.875 -> 3909.6921
$ws.Cells.Item(127, 13).Value = 3434.00002  # M127: 3156.4999 -> 3434.00002
$ws.Cells.Item(127, 14).Value = -13829.6921  # N127: -13806.875 -> -13829.6921

$ws.Cells.Item(129, 8).Value = 1550.0454  # H129: 1255.909 -> 1550.0454
$ws.Cells.Item(129, 9).Value = 591.3333  # I129: 613 -> 591.3333
$ws.Cells.Item(129, 10).Value = 1909.5625  # J129: 1338.3334 -> 1909.5625
$ws.Cells.Item(129, 11).Value = 1773.9999  # K129: 1839 -> 1773.9999
$ws.Cells.Item(129, 12).Value = 5728.6875  # L129: 4015.0002 -> 5728.6875
$ws.Cells.Item(129, 13).Value = 3226.0001  # M129: 3161 -> 3226.0001
$ws.Cells.Item(129, 14).Value = -15728.6875  # N129: -14015.0002 -> -15728.6875

$ws.Cells.Item(132, 8).Value = 2537.0967  # H132: 2638.5925 -> 2537.0967
$ws.Cells.Item(132, 9).Value = 2631.4814  # I132: 2638.5925 -> 2631.4814
$ws.Cells.Item(132, 10).Value = 1900  # J132: 0 -> 1900
$ws.Cells.Item(132, 11).Value = 7894.4442  # K132: 7915.7775 -> 7894.4442
$ws.Cells.Item(132, 12).Value = 5700  # L132: 0 -> 5700
$ws.Cells.Item(132, 13).Value = -5364.4442  # M132: -5385.7775 -> -5364.4442
$ws.Cells.Item(132, 14).Value = -10760  # N132: None -> -10760

$ws.Cells.Item(136, 8).Value = 57223.332  # H136: 30000 -> 57223.332
$ws.Cells.Item(136, 9).Value = 35000  # I136: 0 -> 35000
$ws.Cells.Item(136, 10).Value = 61668  # J136: 30000 -> 61668
$ws.Cells.Item(136, 11).Value = 35000  # K136: 0 -> 35000
$ws.Cells.Item(136, 12).Value = 61668  # L136: 30000 -> 61668
$ws.Cells.Item(136, 13).Value = -29900  # M136: None -> -29900
$ws.Cells.Item(136, 14).Value = -71868  # N136: -40200 -> -71868

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(35, 8).Value = 13918.25  # H35: 10035.833 -> 13918.25
$ws.Cells.Item(35, 9).Value = 13918.25  # I35: 13918.5 -> 13918.25
$ws.Cells.Item(35, 10).Value = 0  # J35: 2270.5 -> 0
$ws.Cells.Item(35, 11).Value = 13918.25  # K35: 13918.5 -> 13918.25
$ws.Cells.Item(35, 12).Value = 0  # L35: 2270.5 -> 0
$ws.Cells.Item(35, 13).Value = -13512.25  # M35: -13512.5 -> -13512.25
$ws.Cells.Item(35, 14).ClearContents()  # N35: -3082.5 -> (removed)

$ws.Cells.Item(61, 8).Value = 3122.7144  # H61: 2085.28 -> 3122.7144
$ws.Cells.Item(61, 9).Value = 3440.75  # I61: 1941.1111 -> 3440.75
$ws.Cells.Item(61, 10).Value = 2698.6667  # J61: 2456 -> 2698.6667
$ws.Cells.Item(61, 11).Value = 3440.75  # K61: 1941.1111 -> 3440.75
$ws.Cells.Item(61, 12).Value = 2698.6667  # L61: 2456 -> 2698.6667
$ws.Cells.Item(61, 13).Value = -3228.75  # M61: -1729.1111 -> -3228.75
$ws.Cells.Item(61, 14).Value = -3122.6667  # N61: -2880 -> -3122.6667

$ws.Cells.Item(122, 8).Value = 3430.3914  # H122: 3682.1428 -> 3430.3914
$ws.Cells.Item(122, 9).Value = 3216.125  # I122: 3563.1428 -> 3216.125
$ws.Cells.Item(122, 11).Value = 9648.375  # K122: 10689.4284 -> 9648.375
$ws.Cells.Item(122, 13).Value = -7198.375  # M122: -8239.428400000001 -> -7198.375

$ws.Cells.Item(132, 8).Value = 4037.8667  # H132: 4444.396 -> 4037.8667
$ws.Cells.Item(132, 9).Value = 5020.8486  # I132: 5748.1787 -> 5020.8486
$ws.Cells.Item(132, 10).Value = 2836.4443  # J132: 2984.16 -> 2836.4443
$ws.Cells.Item(132, 11).Value = 15062.5458  # K132: 17244.5361 -> 15062.5458
$ws.Cells.Item(132, 12).Value = 8509.332900000001  # L132: 8952.48 -> 8509.332900000001
$ws.Cells.Item(132, 13).Value = -12532.5458  # M132: -14714.5361 -> -12532.5458
$ws.Cells.Item(132, 14).Value = -13569.3329  # N132: -14012.48 -> -13569.3329

$ws.Cells.Item(134, 8).Value = 19473  # H134: 44209.5 -> 19473
$ws.Cells.Item(134, 10).Value = 19473  # J134: 44209.5 -> 19473
$ws.Cells.Item(134, 12).Value = 19473  # L134: 44209.5 -> 19473
$ws.Cells.Item(134, 14).Value = -29613  # N134: -54349.5 -> -29613

$ws.Cells.Item(136, 8).Value = 3122.7144  # H136: 2085.28 -> 3122.7144
$ws.Cells.Item(136, 9).Value = 3440.75  # I136: 1941.1111 -> 3440.75
$ws.Cells.Item(136, 10).Value = 2698.6667  # J136: 2456 -> 2698.6667
$ws.Cells.Item(136, 11).Value = 10322.25  # K136: 5823.3333 -> 10322.25
$ws.Cells.Item(136, 12).Value = 8096.000100000001  # L136: 7368 -> 8096.000100000001
$ws.Cells.Item(136, 13).Value = -7772.25  # M136: -3273.3333 -> -7772.25
$ws.Cells.Item(136, 14).Value = -13196.0001  # N136: -12468 -> -13196.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(62, 8).Value = 30499.5  # H62: 35000 -> 30499.5
$ws.Cells.Item(62, 10).Value = 30499.5  # J62: 35000 -> 30499.5
$ws.Cells.Item(62, 12).Value = 30499.5  # L62: 35000 -> 30499.5
$ws.Cells.Item(62, 14).Value = -31871.5  # N62: -36372 -> -31871.5

$ws.Cells.Item(65, 8).Value = 30499.5  # H65: 35000 -> 30499.5
$ws.Cells.Item(65, 10).Value = 30499.5  # J65: 35000 -> 30499.5
$ws.Cells.Item(65, 12).Value = 91498.5  # L65: 105000 -> 91498.5
$ws.Cells.Item(65, 14).Value = -98362.5  # N65: -111864 -> -98362.5

$ws.Cells.Item(94, 8).Value = 1978.25  # H94: 2120.9285 -> 1978.25
$ws.Cells.Item(94, 9).Value = 1296.2858  # I94: 1423 -> 1296.2858
$ws.Cells.Item(94, 11).Value = 1296.2858  # K94: 1423 -> 1296.2858
$ws.Cells.Item(94, 13).Value = -845.2858000000001  # M94: -972 -> -845.2858000000001

$ws.Cells.Item(107, 8).Value = 15055.308  # H107: 16693.6 -> 15055.308
$ws.Cells.Item(107, 9).Value = 18229.967  # I107: 19487.861 -> 18229.967
$ws.Cells.Item(107, 10).Value = 2753.5  # J107: 3188 -> 2753.5
$ws.Cells.Item(107, 11).Value = 18229.967  # K107: 19487.861 -> 18229.967
$ws.Cells.Item(107, 12).Value = 2753.5  # L107: 3188 -> 2753.5
$ws.Cells.Item(107, 13).Value = -16309.967  # M107: -17567.861 -> -16309.967
$ws.Cells.Item(107, 14).Value = -6593.5  # N107: -7028 -> -6593.5

$ws.Cells.Item(135, 8).Value = 55015  # H135: 75160 -> 55015
$ws.Cells.Item(135, 10).Value = 55015  # J135: 75160 -> 55015
$ws.Cells.Item(135, 12).Value = 55015  # L135: 75160 -> 55015
$ws.Cells.Item(135, 14).Value = -65155  # N135: -85300 -> -65155

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 8775  # H4: 9999.556 -> 8775
$ws.Cells.Item(4, 9).Value = 8162.5  # I4: 0 -> 8162.5
$ws.Cells.Item(4, 10).Value = 10000  # J4: 9999.556 -> 10000
$ws.Cells.Item(4, 11).Value = 8162.5  # K4: 0 -> 8162.5
$ws.Cells.Item(4, 12).Value = 10000  # L4: 9999.556 -> 10000
$ws.Cells.Item(4, 13).Value = -8050.5  # M4: None -> -8050.5
$ws.Cells.Item(4, 14).Value = -10224  # N4: -10223.556 -> -10224

$ws.Cells.Item(58, 8).Value = 1581.3636  # H58: 1725.3846 -> 1581.3636
$ws.Cells.Item(58, 9).Value = 1497.8889  # I58: 1639.4584 -> 1497.8889
$ws.Cells.Item(58, 10).Value = 1713.9412  # J58: 1862.8667 -> 1713.9412
$ws.Cells.Item(58, 11).Value = 1497.8889  # K58: 1639.4584 -> 1497.8889
$ws.Cells.Item(58, 12).Value = 1713.9412  # L58: 1862.8667 -> 1713.9412
$ws.Cells.Item(58, 13).Value = -1294.8889  # M58: -1436.4584 -> -1294.8889
$ws.Cells.Item(58, 14).Value = -2119.9412  # N58: -2268.8667 -> -2119.9412

$ws.Cells.Item(68, 8).Value = 24314  # H68: 28333 -> 24314
$ws.Cells.Item(68, 10).Value = 24314  # J68: 28333 -> 24314
$ws.Cells.Item(68, 12).Value = 24314  # L68: 28333 -> 24314
$ws.Cells.Item(68, 14).Value = -25812  # N68: -29831 -> -25812

$ws.Cells.Item(71, 8).Value = 24314  # H71: 28333 -> 24314
$ws.Cells.Item(71, 10).Value = 24314  # J71: 28333 -> 24314
$ws.Cells.Item(71, 12).Value = 72942  # L71: 84999 -> 72942
$ws.Cells.Item(71, 14).Value = -80430  # N71: -92487 -> -80430

$ws.Cells.Item(97, 8).Value = 39499.5  # H97: 30999.5 -> 39499.5
$ws.Cells.Item(97, 10).Value = 39499.5  # J97: 30999.5 -> 39499.5
$ws.Cells.Item(97, 12).Value = 39499.5  # L97: 30999.5 -> 39499.5
$ws.Cells.Item(97, 14).Value = -41481.5  # N97: -32981.5 -> -41481.5

$ws.Cells.Item(105, 8).Value = 1578.5358  # H105: 1711.625 -> 1578.5358
$ws.Cells.Item(105, 9).Value = 1558.2916  # I105: 1713.95 -> 1558.2916
$ws.Cells.Item(105, 11).Value = 1558.2916  # K105: 1713.95 -> 1558.2916
$ws.Cells.Item(105, 13).Value = 188.7084  # M105: 33.04999999999995 -> 188.7084

$ws.Cells.Item(134, 8).Value = 1701.4324  # H134: 2194.4814 -> 1701.4324
$ws.Cells.Item(134, 9).Value = 1219.3214  # I134: 1582.1666 -> 1219.3214
$ws.Cells.Item(134, 10).Value = 3201.3333  # J134: 3419.111 -> 3201.3333
$ws.Cells.Item(134, 11).Value = 3657.9642  # K134: 4746.4998 -> 3657.9642
$ws.Cells.Item(134, 12).Value = 9603.999899999999  # L134: 10257.333 -> 9603.999899999999
$ws.Cells.Item(134, 13).Value = -1122.9642  # M134: -2211.4998 -> -1122.9642
$ws.Cells.Item(134, 14).Value = -14673.9999  # N134: -15327.333 -> -14673.9999

$ws.Cells.Item(136, 8).Value = 1581.3636  # H136: 1725.3846 -> 1581.3636
$ws.Cells.Item(136, 9).Value = 1497.8889  # I136: 1639.4584 -> 1497.8889
$ws.Cells.Item(136, 10).Value = 1713.9412  # J136: 1862.8667 -> 1713.9412
$ws.Cells.Item(136, 11).Value = 4493.6667  # K136: 4918.3752 -> 4493.6667
$ws.Cells.Item(136, 12).Value = 5141.8236  # L136: 5588.6001 -> 5141.8236
$ws.Cells.Item(136, 13).Value = -1943.6667  # M136: -2368.3752 -> -1943.6667
$ws.Cells.Item(136, 14).Value = -10241.8236  # N136: -10688.6001 -> -10241.8236

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 0  # H7: 177.27777 -> 0
$ws.Cells.Item(7, 9).Value = 0  # I7: 131.61539 -> 0
$ws.Cells.Item(7, 10).Value = 0  # J7: 296 -> 0
$ws.Cells.Item(7, 11).Value = 0  # K7: 394.84617 -> 0
$ws.Cells.Item(7, 12).Value = 0  # L7: 888 -> 0
$ws.Cells.Item(7, 13).ClearContents()  # M7: -282.84617 -> (removed)
$ws.Cells.Item(7, 14).ClearContents()  # N7: -1112 -> (removed)

$ws.Cells.Item(80, 8).Value = 11000  # H80: 9741.5 -> 11000
$ws.Cells.Item(80, 9).Value = 15000  # I80: 12223.75 -> 15000
$ws.Cells.Item(80, 10).Value = 9285.714  # J80: 8500.375 -> 9285.714
$ws.Cells.Item(80, 11).Value = 45000  # K80: 36671.25 -> 45000
$ws.Cells.Item(80, 12).Value = 27857.142  # L80: 25501.125 -> 27857.142
$ws.Cells.Item(80, 13).Value = -44064  # M80: -35735.25 -> -44064
$ws.Cells.Item(80, 14).Value = -29729.142  # N80: -27373.125 -> -29729.142

$ws.Cells.Item(83, 8).Value = 11000  # H83: 9741.5 -> 11000
$ws.Cells.Item(83, 9).Value = 15000  # I83: 12223.75 -> 15000
$ws.Cells.Item(83, 10).Value = 9285.714  # J83: 8500.375 -> 9285.714
$ws.Cells.Item(83, 11).Value = 135000  # K83: 110013.75 -> 135000
$ws.Cells.Item(83, 12).Value = 83571.42600000001  # L83: 76503.375 -> 83571.42600000001
$ws.Cells.Item(83, 13).Value = -130320  # M83: -105333.75 -> -130320
$ws.Cells.Item(83, 14).Value = -92931.42600000001  # N83: -85863.375 -> -92931.42600000001

$ws.Cells.Item(98, 8).Value = 2124.3333  # H98: 558.9286 -> 2124.3333
$ws.Cells.Item(98, 9).Value = 2811.5  # I98: 453.16666 -> 2811.5
$ws.Cells.Item(98, 10).Value = 750  # J98: 638.25 -> 750
$ws.Cells.Item(98, 11).Value = 8434.5  # K98: 1359.49998 -> 8434.5
$ws.Cells.Item(98, 12).Value = 2250  # L98: 1914.75 -> 2250
$ws.Cells.Item(98, 13).Value = -6936.5  # M98: 138.5000199999999 -> -6936.5
$ws.Cells.Item(98, 14).Value = -5246  # N98: -4910.75 -> -5246

$ws.Cells.Item(131, 8).Value = 785.62  # H131: 10640712 -> 785.62
$ws.Cells.Item(131, 9).Value = 485.8  # I131: 17071.5 -> 485.8
$ws.Cells.Item(131, 10).Value = 801.4  # J131: 11365051 -> 801.4
$ws.Cells.Item(131, 11).Value = 1457.4  # K131: 51214.5 -> 1457.4
$ws.Cells.Item(131, 12).Value = 2404.2  # L131: 34095153 -> 2404.2
$ws.Cells.Item(131, 13).Value = 3582.6  # M131: -46174.5 -> 3582.6
$ws.Cells.Item(131, 14).Value = -12484.2  # N131: -34105233 -> -12484.2

$ws.Cells.Item(133, 8).Value = 5272.857  # H133: 5807.273 -> 5272.857
$ws.Cells.Item(133, 9).Value = 1987  # I133: 1990.5 -> 1987
$ws.Cells.Item(133, 10).Value = 6169  # J133: 6655.4443 -> 6169
$ws.Cells.Item(133, 11).Value = 5961  # K133: 5971.5 -> 5961
$ws.Cells.Item(133, 12).Value = 18507  # L133: 19966.3329 -> 18507
$ws.Cells.Item(133, 13).Value = -901  # M133: -911.5 -> -901
$ws.Cells.Item(133, 14).Value = -28627  # N133: -30086.3329 -> -28627

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 40000  # H57: 0 -> 40000
$ws.Cells.Item(57, 10).Value = 40000  # J57: 0 -> 40000
$ws.Cells.Item(57, 12).Value = 40000  # L57: 0 -> 40000
$ws.Cells.Item(57, 14).Value = -41640  # N57: None -> -41640

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1286.4546  # H46: 1428.4286 -> 1286.4546
$ws.Cells.Item(46, 9).Value = 1100.2  # I46: 1125 -> 1100.2
$ws.Cells.Item(46, 10).Value = 1441.6666  # J46: 1549.8 -> 1441.6666
$ws.Cells.Item(46, 11).Value = 1100.2  # K46: 1125 -> 1100.2
$ws.Cells.Item(46, 12).Value = 1441.6666  # L46: 1549.8 -> 1441.6666
$ws.Cells.Item(46, 13).Value = -912.2  # M46: -937 -> -912.2
$ws.Cells.Item(46, 14).Value = -1817.6666  # N46: -1925.8 -> -1817.6666

$ws.Cells.Item(62, 8).Value = 35000  # H62: 0 -> 35000
$ws.Cells.Item(62, 9).Value = 35000  # I62: 0 -> 35000
$ws.Cells.Item(62, 11).Value = 35000  # K62: 0 -> 35000
$ws.Cells.Item(62, 13).Value = -34376  # M62: None -> -34376

$ws.Cells.Item(65, 8).Value = 35000  # H65: 0 -> 35000
$ws.Cells.Item(65, 9).Value = 35000  # I65: 0 -> 35000
$ws.Cells.Item(65, 11).Value = 105000  # K65: 0 -> 105000
$ws.Cells.Item(65, 13).Value = -101880  # M65: None -> -101880

$ws.Cells.Item(122, 8).Value = 20005370  # H122: 14290193 -> 20005370
$ws.Cells.Item(122, 9).Value = 4814.2856  # I122: 3882.5454 -> 4814.2856
$ws.Cells.Item(122, 11).Value = 14442.8568  # K122: 11647.6362 -> 14442.8568
$ws.Cells.Item(122, 13).Value = -11992.8568  # M122: -9197.636200000001 -> -11992.8568

$ws.Cells.Item(135, 8).Value = 54995  # H135: 0 -> 54995
$ws.Cells.Item(135, 10).Value = 54995  # J135: 0 -> 54995
$ws.Cells.Item(135, 12).Value = 54995  # L135: 0 -> 54995
$ws.Cells.Item(135, 14).Value = -65135  # N135: None -> -65135

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(11, 8).Value = 49999.5  # H11: 20000 -> 49999.5
$ws.Cells.Item(11, 10).Value = 49999  # J11: 10000 -> 49999
$ws.Cells.Item(11, 12).Value = 49999  # L11: 10000 -> 49999
$ws.Cells.Item(11, 14).Value = -50283  # N11: -10284 -> -50283

$ws.Cells.Item(15, 8).Value = 55002.5  # H15: 55000 -> 55002.5
$ws.Cells.Item(15, 10).Value = 55002.5  # J15: 55000 -> 55002.5
$ws.Cells.Item(15, 12).Value = 55002.5  # L15: 55000 -> 55002.5
$ws.Cells.Item(15, 14).Value = -55578.5  # N15: -55576 -> -55578.5

$ws.Cells.Item(20, 8).Value = 38333.332  # H20: 32505.5 -> 38333.332
$ws.Cells.Item(20, 10).Value = 15000  # J20: 15011 -> 15000
$ws.Cells.Item(20, 12).Value = 15000  # L20: 15011 -> 15000
$ws.Cells.Item(20, 14).Value = -15480  # N20: -15491 -> -15480

$ws.Cells.Item(21, 8).Value = 19469.166  # H21: 30750 -> 19469.166
$ws.Cells.Item(21, 9).Value = 19469.166  # I21: 30750 -> 19469.166
$ws.Cells.Item(21, 11).Value = 19469.166  # K21: 30750 -> 19469.166
$ws.Cells.Item(21, 13).Value = -19234.166  # M21: -30515 -> -19234.166

$ws.Cells.Item(35, 8).Value = 19469.166  # H35: 30750 -> 19469.166
$ws.Cells.Item(35, 9).Value = 19469.166  # I35: 30750 -> 19469.166
$ws.Cells.Item(35, 11).Value = 19469.166  # K35: 30750 -> 19469.166
$ws.Cells.Item(35, 13).Value = -19179.166  # M35: -30460 -> -19179.166

$ws.Cells.Item(81, 8).Value = 83054.625  # H81: 63627.75 -> 83054.625
$ws.Cells.Item(81, 9).Value = 72633.86  # I81: 63627.75 -> 72633.86
$ws.Cells.Item(81, 10).Value = 156000  # J81: 0 -> 156000
$ws.Cells.Item(81, 11).Value = 145267.72  # K81: 127255.5 -> 145267.72
$ws.Cells.Item(81, 12).Value = 312000  # L81: 0 -> 312000
$ws.Cells.Item(81, 13).Value = -144206.72  # M81: -126194.5 -> -144206.72
$ws.Cells.Item(81, 14).Value = -314122  # N81: None -> -314122

$ws.Cells.Item(84, 8).Value = 83054.625  # H84: 63627.75 -> 83054.625
$ws.Cells.Item(84, 9).Value = 72633.86  # I84: 63627.75 -> 72633.86
$ws.Cells.Item(84, 10).Value = 156000  # J84: 0 -> 156000
$ws.Cells.Item(84, 11).Value = 726338.6  # K84: 636277.5 -> 726338.6
$ws.Cells.Item(84, 12).Value = 1560000  # L84: 0 -> 1560000
$ws.Cells.Item(84, 13).Value = -721034.6  # M84: -630973.5 -> -721034.6
$ws.Cells.Item(84, 14).Value = -1570608  # N84: None -> -1570608

$ws.Cells.Item(135, 8).Value = 57905  # H135: 54905 -> 57905
$ws.Cells.Item(135, 10).Value = 57905  # J135: 54905 -> 57905
$ws.Cells.Item(135, 12).Value = 57905  # L135: 54905 -> 57905
$ws.Cells.Item(135, 14).Value = -68045  # N135: -65045 -> -68045

$ws.Cells.Item(136, 8).Value = 2121.5186  # H136: 2184.6924 -> 2121.5186
$ws.Cells.Item(136, 9).Value = 1832.3182  # I136: 1896.762 -> 1832.3182
$ws.Cells.Item(136, 11).Value = 5496.9546  # K136: 5690.286 -> 5496.9546
$ws.Cells.Item(136, 13).Value = -2946.9546  # M136: -3140.286 -> -2946.9546
